# "Replaced pi with h": the hadron label "pi+" used in the dataset becomes "h+".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq "pi+") {
            $cell.Value = "h+"
        }
    }
}
